$p = $ppt.ActivePresentation
$s = $p.Slides.Item(4)
$sh = $s.Shapes.Item(2)
$tr = $sh.TextFrame.TextRange

# Insert the new paragraph "不可包含任何合并单元格" right before the
# "数据列数需 >= 231" paragraph (currently paragraph 2).
$targetPara = $tr.Paragraphs(2, 1)
[void]$targetPara.InsertBefore("不可包含任何合并单元格`r")

# Paragraphs 1-7 (the sz=2400 text block, now including the newly
# inserted paragraph) shrink from 24pt to 20pt.
$firstPara = $tr.Paragraphs(1, 1)
$lastBigPara = $tr.Paragraphs(7, 1)
$bigStart = $firstPara.Start
$bigEnd = $lastBigPara.Start + $lastBigPara.Length - 1
$bigRange = $tr.Characters($bigStart, $bigEnd - $bigStart + 1)
$bigRange.Font.Size = 20

# Paragraphs 8-10 (the lvl=1 sz=1400 bullet sub-points) shrink from
# 14pt to 12pt.
$firstSmallPara = $tr.Paragraphs(8, 1)
$lastSmallPara = $tr.Paragraphs(10, 1)
$smallStart = $firstSmallPara.Start
$smallEnd = $lastSmallPara.Start + $lastSmallPara.Length - 1
$smallRange = $tr.Characters($smallStart, $smallEnd - $smallStart + 1)
$smallRange.Font.Size = 12
